$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")
$ws.Range("H8").Value = 0.396
